# Auto-generated edit script
# Applies updated market-price / profit figures to several leve sheets
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H6").Value = 12500333
$ws.Range("I6").Value = 14285809
$ws.Range("J6").Value = 2001
$ws.Range("K6").Value = 42857427
$ws.Range("L6").Value = 6003
$ws.Range("M6").Value = -42857315
$ws.Range("N6").Value = -6227
$ws.Range("H8").Value = 66.166664
$ws.Range("I8").Value = 66.166664
$ws.Range("K8").Value = 198.499992
$ws.Range("M8").Value = -59.49999199999999
$ws.Range("H13").Value = 28405
$ws.Range("J13").Value = 28405
$ws.Range("L13").Value = 28405
$ws.Range("N13").Value = -28743
$ws.Range("H21").Value = 20192
$ws.Range("I21").Value = 19307.2
$ws.Range("J21").Value = 21666.666
$ws.Range("K21").Value = 19307.2
$ws.Range("L21").Value = 21666.666
$ws.Range("M21").Value = -18839.2
$ws.Range("N21").Value = -22602.666
$ws.Range("H23").Value = 20192
$ws.Range("I23").Value = 19307.2
$ws.Range("J23").Value = 21666.666
$ws.Range("K23").Value = 19307.2
$ws.Range("L23").Value = 21666.666
$ws.Range("M23").Value = -19073.2
$ws.Range("N23").Value = -22134.666
$ws.Range("H69").Value = 13238.8
$ws.Range("I69").Value = 1400
$ws.Range("J69").Value = 14554.223
$ws.Range("K69").Value = 4200
$ws.Range("L69").Value = 43662.669
$ws.Range("M69").Value = -3326
$ws.Range("N69").Value = -45410.669
$ws.Range("H72").Value = 13238.8
$ws.Range("I72").Value = 1400
$ws.Range("J72").Value = 14554.223
$ws.Range("K72").Value = 12600
$ws.Range("L72").Value = 130988.007
$ws.Range("M72").Value = -8232
$ws.Range("N72").Value = -139724.007
$ws.Range("H76").Value = 3066
$ws.Range("I76").Value = 3041.8604
$ws.Range("J76").Value = 3214.2856
$ws.Range("K76").Value = 3041.8604
$ws.Range("L76").Value = 3214.2856
$ws.Range("M76").Value = -2726.8604
$ws.Range("N76").Value = -3844.2856
$ws.Range("H79").Value = 3066
$ws.Range("I79").Value = 3041.8604
$ws.Range("J79").Value = 3214.2856
$ws.Range("K79").Value = 3041.8604
$ws.Range("L79").Value = 3214.2856
$ws.Range("M79").Value = -1949.8604
$ws.Range("N79").Value = -5398.2856
$ws.Range("H80").Value = 4090.78
$ws.Range("I80").Value = 3245.2
$ws.Range("K80").Value = 9735.599999999999
$ws.Range("M80").Value = -8737.599999999999
$ws.Range("H83").Value = 4090.78
$ws.Range("I83").Value = 3245.2
$ws.Range("K83").Value = 29206.8
$ws.Range("M83").Value = -24214.8
$ws.Range("H88").Value = 2748.8076
$ws.Range("I88").Value = 1842.375
$ws.Range("J88").Value = 3151.6667
$ws.Range("K88").Value = 1842.375
$ws.Range("L88").Value = 3151.6667
$ws.Range("M88").Value = -1436.375
$ws.Range("N88").Value = -3963.6667
$ws.Range("H91").Value = 2748.8076
$ws.Range("I91").Value = 1842.375
$ws.Range("J91").Value = 3151.6667
$ws.Range("K91").Value = 1842.375
$ws.Range("L91").Value = 3151.6667
$ws.Range("M91").Value = -438.375
$ws.Range("N91").Value = -5959.6667
$ws.Range("H121").Value = 1162.75
$ws.Range("I121").Value = 201
$ws.Range("J121").Value = 1483.3334
$ws.Range("K121").Value = 603
$ws.Range("L121").Value = 4450.0002
$ws.Range("M121").Value = 1144
$ws.Range("N121").Value = -7944.0002
$ws.Range("H132").Value = 2037.9166
$ws.Range("I132").Value = 1327.9032
$ws.Range("J132").Value = 6440
$ws.Range("K132").Value = 3983.7096
$ws.Range("L132").Value = 19320
$ws.Range("M132").Value = -1453.7096
$ws.Range("N132").Value = -24380

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H10").Value = 70005
$ws.Range("I10").Value = 0
$ws.Range("J10").Value = 70005
$ws.Range("K10").Value = 0
$ws.Range("L10").Value = 70005
$ws.Range("M10").ClearContents()
$ws.Range("N10").Value = -70345
$ws.Range("H32").Value = 3287.5
$ws.Range("I32").Value = 2911.3147
$ws.Range("J32").Value = 5826.75
$ws.Range("K32").Value = 2911.3147
$ws.Range("L32").Value = 5826.75
$ws.Range("M32").Value = -2624.3147
$ws.Range("N32").Value = -6400.75
$ws.Range("H88").Value = 2183.8333
$ws.Range("I88").Value = 1923.5555
$ws.Range("J88").Value = 2340
$ws.Range("K88").Value = 1923.5555
$ws.Range("L88").Value = 2340
$ws.Range("M88").Value = -1517.5555
$ws.Range("N88").Value = -3152
$ws.Range("H91").Value = 2183.8333
$ws.Range("I91").Value = 1923.5555
$ws.Range("J91").Value = 2340
$ws.Range("K91").Value = 1923.5555
$ws.Range("L91").Value = 2340
$ws.Range("M91").Value = -519.5554999999999
$ws.Range("N91").Value = -5148

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H19").Value = 21466.924
$ws.Range("J19").Value = 21466.924
$ws.Range("L19").Value = 21466.924
$ws.Range("N19").Value = -21812.924

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H23").Value = 11066.667
$ws.Range("I23").Value = 5200
$ws.Range("J23").Value = 14000
$ws.Range("K23").Value = 5200
$ws.Range("L23").Value = 14000
$ws.Range("M23").Value = -4960
$ws.Range("N23").Value = -14480
$ws.Range("H27").Value = 11066.667
$ws.Range("I27").Value = 5200
$ws.Range("J27").Value = 14000
$ws.Range("K27").Value = 5200
$ws.Range("L27").Value = 14000
$ws.Range("M27").Value = -5008
$ws.Range("N27").Value = -14384
$ws.Range("H31").Value = 21740516
$ws.Range("I31").Value = 1012.1539
$ws.Range("J31").Value = 50001870
$ws.Range("K31").Value = 1012.1539
$ws.Range("L31").Value = 50001870
$ws.Range("M31").Value = -717.1539
$ws.Range("N31").Value = -50002460
$ws.Range("H34").Value = 21740516
$ws.Range("I34").Value = 1012.1539
$ws.Range("J34").Value = 50001870
$ws.Range("K34").Value = 1012.1539
$ws.Range("L34").Value = 50001870
$ws.Range("M34").Value = -810.1539
$ws.Range("N34").Value = -50002274
$ws.Range("H36").Value = 14333.167
$ws.Range("J36").Value = 14333.167
$ws.Range("L36").Value = 14333.167
$ws.Range("N36").Value = -15109.167
$ws.Range("H40").Value = 14333.167
$ws.Range("J40").Value = 14333.167
$ws.Range("L40").Value = 14333.167
$ws.Range("N40").Value = -14653.167
$ws.Range("H107").Value = 1106.0769
$ws.Range("I107").Value = 422.58334
$ws.Range("J107").Value = 1691.9286
$ws.Range("K107").Value = 422.58334
$ws.Range("L107").Value = 1691.9286
$ws.Range("M107").Value = 1497.41666
$ws.Range("N107").Value = -5531.9286
$ws.Range("H137").Value = 33748
$ws.Range("J137").Value = 33748
$ws.Range("L137").Value = 33748
$ws.Range("N137").Value = -43948

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H7").Value = 8695764
$ws.Range("I7").Value = 68.875
$ws.Range("J7").Value = 13333468
$ws.Range("K7").Value = 206.625
$ws.Range("L7").Value = 40000404
$ws.Range("M7").Value = -94.625
$ws.Range("N7").Value = -40000628
$ws.Range("H12").Value = 234.45833
$ws.Range("I12").Value = 257.77777
$ws.Range("J12").Value = 220.46666
$ws.Range("K12").Value = 773.33331
$ws.Range("L12").Value = 661.3999799999999
$ws.Range("M12").Value = -600.33331
$ws.Range("N12").Value = -1007.39998
$ws.Range("H69").Value = 3000
$ws.Range("H72").Value = 3000
$ws.Range("H92").Value = 3333830.8
$ws.Range("I92").Value = 746
$ws.Range("J92").Value = 10000000
$ws.Range("K92").Value = 2238
$ws.Range("L92").Value = 30000000
$ws.Range("M92").Value = -990
$ws.Range("N92").Value = -30002496
$ws.Range("H97").Value = 554.1957
$ws.Range("I97").Value = 282.72223
$ws.Range("J97").Value = 728.7143
$ws.Range("K97").Value = 848.16669
$ws.Range("L97").Value = 2186.1429
$ws.Range("M97").Value = -352.16669
$ws.Range("N97").Value = -3178.1429
$ws.Range("H113").Value = 614.2593000000001
$ws.Range("I113").Value = 622.7059
$ws.Range("K113").Value = 1868.1177
$ws.Range("M113").Value = 301.8822999999998
$ws.Range("H117").Value = 297.66666
$ws.Range("J117").Value = 650
$ws.Range("L117").Value = 1950
$ws.Range("N117").Value = -8834
$ws.Range("H121").Value = 1167.1428
$ws.Range("I121").Value = 300
$ws.Range("K121").Value = 900
$ws.Range("M121").Value = 410
$ws.Range("H129").Value = 2223.6428
$ws.Range("I129").Value = 2123.75
$ws.Range("J129").Value = 2263.6
$ws.Range("K129").Value = 6371.25
$ws.Range("L129").Value = 6790.799999999999
$ws.Range("M129").Value = -1371.25
$ws.Range("N129").Value = -16790.8
$ws.Range("H131").Value = 1097.2
$ws.Range("J131").Value = 1147.0952
$ws.Range("L131").Value = 3441.2856
$ws.Range("N131").Value = -13521.2856
$ws.Range("H132").Value = 3046
$ws.Range("I132").Value = 2909
$ws.Range("J132").Value = 3320
$ws.Range("K132").Value = 26181
$ws.Range("L132").Value = 29880
$ws.Range("M132").Value = -23651
$ws.Range("N132").Value = -34940

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 1616.6666
$ws.Range("I16").Value = 1616.6666
$ws.Range("K16").Value = 1616.6666
$ws.Range("M16").Value = -1446.6666
$ws.Range("H25").Value = 60000
$ws.Range("J25").Value = 60000
$ws.Range("L25").Value = 60000
$ws.Range("N25").Value = -60460
$ws.Range("H63").Value = 0
$ws.Range("J63").Value = 0
$ws.Range("L63").Value = 0
$ws.Range("N63").ClearContents()
$ws.Range("H66").Value = 0
$ws.Range("J66").Value = 0
$ws.Range("L66").Value = 0
$ws.Range("N66").ClearContents()
$ws.Range("H82").Value = 2533.8462
$ws.Range("I82").Value = 0
$ws.Range("J82").Value = 2533.8462
$ws.Range("K82").Value = 0
$ws.Range("L82").Value = 2533.8462
$ws.Range("M82").ClearContents()
$ws.Range("N82").Value = -3255.8462
$ws.Range("H85").Value = 2533.8462
$ws.Range("I85").Value = 0
$ws.Range("J85").Value = 2533.8462
$ws.Range("K85").Value = 0
$ws.Range("L85").Value = 2533.8462
$ws.Range("M85").ClearContents()
$ws.Range("N85").Value = -5029.8462
$ws.Range("H132").Value = 3494.7368
$ws.Range("I132").Value = 3246.1667
$ws.Range("J132").Value = 3920.8572
$ws.Range("K132").Value = 9738.500100000001
$ws.Range("L132").Value = 11762.5716
$ws.Range("M132").Value = -7208.500100000001
$ws.Range("N132").Value = -16822.5716

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H13").Value = 17150
$ws.Range("I13").Value = 825
$ws.Range("K13").Value = 825
$ws.Range("M13").Value = -685

